$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 81: Course, Hours, and Notes for the JS101 entry
$ws.Range("B81").Value = "JS101: Programming Foundations with JavaScript"
$ws.Range("C81").Value = 0.5
$ws.Range("D81").Value = "Finish 1 small problem"

# Move the active selection to C82, matching the saved cursor position
$ws.Range("C82").Select()
